# Adds a new "CD" column (one further quarterly/period data point) to the
# real-time GVA series sheet, plus one brand-new trailing row (row 140) that
# only has the new date (col A) and the new CD value - mirroring the sparse
# "revision triangle" layout already used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: CD1 is a date, styled like the rest of row 1 (CC1 etc). ---
# Copy CC1's style (number format/font/border/alignment) onto CD1, then set
# its own value so we don't disturb CC1.
$ws.Range("CC1").Copy()
$ws.Range("CD1").PasteSpecial(-4122)
$ws.Range("CD1").Value = 45986

# --- New CD values for data rows 2..140 (index == row number; 0/1 unused). ---
$cdValues = @(
    $null, $null, "70.45999999999999", "69.75", "69.66", "70.61", "71.77", "71.27", "71.13", "70.72",
    "70.31", "70.23999999999999", "70.47", "70.41", "71.04000000000001", "71.73", "72.15000000000001", "73.04000000000001", "72.70999999999999", "73.38",
    "73.56999999999999", "73.73", "72.65000000000001", "74.05", "74.08", "75.08", "74.45999999999999", "75.40000000000001", "75.73", "76.65000000000001",
    "77.09", "76.73999999999999", "77.14", "77.34", "77.88", "77.86", "78.84", "79.18000000000001", "80.51000000000001", "81.11",
    "81.45999999999999", "81.20999999999999", "82.97", "82.58", "82.81", "82.73", "82.40000000000001", "82.72", "83.05", "83.22",
    "81.88", "81.86", "82.89", "82.98", "82.79000000000001", "83.69", "83.29000000000001", "83.28", "83.45999999999999", "83.83",
    "84.41", "84.93000000000001", "85.81", "87.29000000000001", "88.15000000000001", "89.33", "89.59", "90.48999999999999", "90.94", "91.7",
    "92.17", "91.97", "91.33", "89.67", "85.09", "85.27", "85.81999999999999", "86.45", "87.81999999999999", "89.25",
    "90.12", "90.65000000000001", "92.45", "92.62", "93.23999999999999", "93.23999999999999", "93.41", "93.59999999999999", "93.77", "93.39",
    "93.05", "93.97", "94.48999999999999", "94.92", "95.98", "95.64", "96.25", "96.95999999999999", "96.59", "97.23999999999999",
    "97.62", "97.87", "98.98", "99.38", "99.58", "99.98999999999999", "101.22", "102.19", "103.26", "104.13",
    "103.67", "104.2", "103.65", "104.35", "104.65", "104.54", "104.93", "104.77", "102.95", "93.06",
    "101.11", "102", "101.54", "103.49", "104.79", "105.3", "105.91", "105.76", "106.36", "106.44",
    "106.02", "105.27", "105.45", "106.4", "105.59", "105.44", "105.04", "104.75", "105.29", "105.19",
    "105.26"
)

for ($r = 2; $r -le 139; $r++) {
    $v = $cdValues[$r]
    if ($v -ne $null) {
        $ws.Cells.Item($r, 82).Value = [double]$v
    }
}

# --- Brand new row 140: only A140 (date) and CD140 (value) are populated, ---
# --- same sparse pattern as every other row in this sheet. ---
$ws.Range("A139").Copy()
$ws.Range("A140").PasteSpecial(-4122)
$ws.Range("A140").Value = 45884
$ws.Cells.Item(140, 82).Value = [double]$cdValues[140]
